$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new blank column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").EntireColumn.Insert()
$ws.Columns("N").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active sheet/tab, with R8 selected
$ws.Activate()
$ws.Range("R8").Select()
